$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '29.230.76'
Set-TextValue "E2" '  +0.43%  '

# Row 3
Set-TextValue "D3" '1.856.59'
Set-TextValue "E3" '  +0.42%  '

# Row 4
Set-TextValue "E4" '  -0.07%  '

# Row 5
Set-TextValue "D5" '0.7054'
Set-TextValue "E5" '  +1.86%  '

# Row 6
Set-TextValue "D6" '238.25'
Set-TextValue "E6" '  +0.16%  '

# Row 7
Set-TextValue "E7" '  -0.04%  '

# Row 8
Set-TextValue "D8" '0.08007'
Set-TextValue "E8" '  +4.04%  '

# Row 9
Set-TextValue "D9" '0.3022'
Set-TextValue "E9" '  -0.28%  '

# Row 10
Set-TextValue "D10" '23.43'
Set-TextValue "E10" '  +0.94%  '

# Row 11
Set-TextValue "D11" '0.08178'
Set-TextValue "E11" '  +0.82%  '

# Row 12
Set-TextValue "D12" '5.183'
Set-TextValue "E12" '  -0.44%  '

# Row 13
Set-TextValue "B13" 'WrappedEther'
Set-TextValue "C13" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D13" '1.807.66'
Set-TextValue "E13" '  -3.08%  '

# Row 14
Set-TextValue "B14" 'Polygon'
Set-TextValue "C14" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D14" '0.7013'
Set-TextValue "E14" '  -3.14%  '

# Row 15
Set-TextValue "D15" '89.51'
Set-TextValue "E15" '  +0.66%  '

# Row 16
Set-TextValue "D16" '29.059.68'
Set-TextValue "E16" '  -0.12%  '

# Row 17
Set-TextValue "D17" '5.799'
Set-TextValue "E17" '  +1.00%  '

# Row 18
Set-TextValue "D18" '0.000007884'
Set-TextValue "E18" '  +1.44%  '

# Row 19
Set-TextValue "D19" '13.21'
Set-TextValue "E19" '  +0.22%  '

# Row 20
Set-TextValue "D20" '237.24'
Set-TextValue "E20" '  +0.85%  '

# Row 21
Set-TextValue "D21" '0.9992'
Set-TextValue "E21" '  -0.06%  '

# Row 22
Set-TextValue "D22" '1.001'
Set-TextValue "E22" '  -0.07%  '

# Row 23
Set-TextValue "D23" '2.009.78'
Set-TextValue "E23" '  -4.09%  '

# Row 24
Set-TextValue "D24" '7.441'
Set-TextValue "E24" '  -1.92%  '

# Row 25
Set-TextValue "E25" '  +1.16%  '

# Row 26
Set-TextValue "D26" '8.892'
Set-TextValue "E26" '  -0.82%  '

# Row 27
Set-TextValue "D27" '0.1431'
Set-TextValue "E27" '  +0.02%  '

# Row 28
Set-TextValue "E28" '  +0.52%  '

# Row 29
Set-TextValue "D29" '1.916'
Set-TextValue "E29" '  -3.18%  '

# Row 30
Set-TextValue "D30" '1.419'
Set-TextValue "E30" '  +1.47%  '

# Row 31
Set-TextValue "D31" '1.477'
Set-TextValue "E31" '  -0.45%  '

# Row 32
Set-TextValue "D32" '4.355'
Set-TextValue "E32" '  -2.92%  '

# Row 33
Set-TextValue "D33" '4.020'
Set-TextValue "E33" '  +0.13%  '

# Row 34
Set-TextValue "D34" '0.05188'
Set-TextValue "E34" '  -0.57%  '

# Row 35
Set-TextValue "D35" '1.158'
Set-TextValue "E35" '  -2.08%  '

# Row 36
Set-TextValue "D36" '0.7123'
Set-TextValue "E36" '  +2.02%  '

# Row 37
Set-TextValue "D37" '0.9988'
Set-TextValue "E37" '  -2.30%  '

# Row 38
Set-TextValue "D38" '2.648'
Set-TextValue "E38" '  -0.40%  '

# Row 39
Set-TextValue "D39" '0.01854'
Set-TextValue "E39" '  +0.28%  '

# Row 40
Set-TextValue "D40" '2.719'
Set-TextValue "E40" '  +1.51%  '

# Row 41
Set-TextValue "D41" '0.9342'
Set-TextValue "E41" '  +2.14%  '

# Row 42
Set-TextValue "D42" '1.133.60'
Set-TextValue "E42" '  +5.02%  '

# Row 43
Set-TextValue "D43" '5.914'
Set-TextValue "E43" '  -1.36%  '

# Row 44
Set-TextValue "D44" '0.4248'
Set-TextValue "E44" '  -0.06%  '

# Row 45
Set-TextValue "D45" '70.10'
Set-TextValue "E45" '  -0.36%  '

# Row 46
Set-TextValue "D46" '1.001'
Set-TextValue "E46" '  +0.01%  '

# Row 47
Set-TextValue "D47" '102.43'

# Row 48
Set-TextValue "D48" '0.5330'
Set-TextValue "E48" '  -4.49%  '

# Row 49
Set-TextValue "D49" '1.756'
Set-TextValue "E49" '  -0.54%  '

# Row 50
Set-TextValue "D50" '9.177'
Set-TextValue "E50" '  +0.66%  '

# Row 51
Set-TextValue "D51" '6.940'
Set-TextValue "E51" '  -0.67%  '
